$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TemperatureRelay")

# Rename the "IsAvailable" header to "Availability"
$ws.Range("D1").Value = "Availability"

# Change the boolean "TRUE" availability flags into numeric 0 for all data rows
$ws.Range("D2:D11").Value = 0

# Update the selected cell to reflect the author's edit location
$ws.Range("M17").Select()
